$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 1409
$ws.Range("F4").Value = 13598
$ws.Range("F5").Value = 793
$ws.Range("F7").Value = 50
$ws.Range("G7").Value = "不可售"
$ws.Range("F9").Value = 25226
$ws.Range("F10").Value = 556
$ws.Range("F11").Value = 235
$ws.Range("F12").Value = 552
$ws.Range("F13").Value = 152
$ws.Range("F14").Value = 401
$ws.Range("F15").Value = 232
$ws.Range("F16").Value = 335
$ws.Range("F17").Value = 188
$ws.Range("F18").Value = 161
$ws.Range("F19").Value = 39
$ws.Range("F20").Value = 258
$ws.Range("F22").Value = 37
$ws.Range("F23").Value = 1427
$ws.Range("F24").Value = 121
$ws.Range("F27").Value = 111

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 4495
$ws.Range("F3").Value = 214
$ws.Range("F6").Value = 48
$ws.Range("F8").Value = 99
$ws.Range("F9").Value = 99
$ws.Range("F17").Value = 23

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 912
$ws.Range("F3").Value = 4811
$ws.Range("F4").Value = 153

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 912
$ws.Range("F4").Value = 1409
$ws.Range("F5").Value = 13598
$ws.Range("F6").Value = 793
$ws.Range("F7").Value = 4811
$ws.Range("F9").Value = 50
$ws.Range("G9").Value = "不可售"
$ws.Range("F11").Value = 153
$ws.Range("F12").Value = 25226
$ws.Range("F13").Value = 556
$ws.Range("F14").Value = 4495
$ws.Range("F15").Value = 235
$ws.Range("F16").Value = 214
$ws.Range("F17").Value = 214
$ws.Range("F18").Value = 552
$ws.Range("F21").Value = 152
$ws.Range("F22").Value = 48
$ws.Range("F23").Value = 48
$ws.Range("F25").Value = 99
$ws.Range("F26").Value = 99
$ws.Range("F28").Value = 401
$ws.Range("F30").Value = 232
$ws.Range("F31").Value = 335
$ws.Range("F32").Value = 188
$ws.Range("F33").Value = 161
$ws.Range("F34").Value = 39
$ws.Range("F36").Value = 258
$ws.Range("F40").Value = 37
$ws.Range("F42").Value = 1427
$ws.Range("F43").Value = 121
$ws.Range("F47").Value = 111
$ws.Range("F48").Value = 23
